$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New schedule rows for group INT-B2-B (rows 52-91), mirroring the existing
# INT-B2-A block's layout/formatting but with its own subject/date sequence.
# Tuple layout: (rowNumber, Year, Group, Subject, Session, Date, StartTime, DurationMinutes)
$newRows = @(
  @(52, "Year 4", "INT-B2-B", "cardiology", "1", "06/12/2025", "08:30:00", 720),
  @(53, "Year 4", "INT-B2-B", "cardiology", "2", "07/12/2025", "08:30:00", 720),
  @(54, "Year 4", "INT-B2-B", "cardiology", "3", "08/12/2025", "08:30:00", 720),
  @(55, "Year 4", "INT-B2-B", "cardiology", "4", "09/12/2025", "08:30:00", 720),
  @(56, "Year 4", "INT-B2-B", "cardiology", "5", "10/12/2025", "08:30:00", 720),
  @(57, "Year 4", "INT-B2-B", "cardiology", "6", "13/12/2025", "08:30:00", 720),
  @(58, "Year 4", "INT-B2-B", "cardiology", "7", "14/12/2025", "08:30:00", 720),
  @(59, "Year 4", "INT-B2-B", "cardiology", "8", "15/12/2025", "08:30:00", 720),
  @(60, "Year 4", "INT-B2-B", "cardiology", "9", "16/12/2025", "08:30:00", 720),
  @(61, "Year 4", "INT-B2-B", "cardiology", "10", "17/12/2025", "08:30:00", 720),
  @(62, "Year 4", "INT-B2-B", "chest", "1", "21/12/2025", "08:30:00", 720),
  @(63, "Year 4", "INT-B2-B", "chest", "2", "22/12/2025", "08:30:00", 720),
  @(64, "Year 4", "INT-B2-B", "chest", "3", "23/12/2025", "08:30:00", 720),
  @(65, "Year 4", "INT-B2-B", "chest", "4", "24/12/2025", "08:30:00", 720),
  @(66, "Year 4", "INT-B2-B", "chest", "5", "25/12/2025", "08:30:00", 720),
  @(67, "Year 4", "INT-B2-B", "chest", "6", "28/12/2025", "08:30:00", 720),
  @(68, "Year 4", "INT-B2-B", "chest", "7", "29/12/2025", "08:30:00", 720),
  @(69, "Year 4", "INT-B2-B", "chest", "8", "30/12/2025", "08:30:00", 720),
  @(70, "Year 4", "INT-B2-B", "chest", "9", "31/12/2025", "08:30:00", 720),
  @(71, "Year 4", "INT-B2-B", "chest", "10", "01/01/2026", "08:30:00", 720),
  @(72, "Year 4", "INT-B2-B", "dermatology", "1", "04/01/2026", "08:30:00", 720),
  @(73, "Year 4", "INT-B2-B", "dermatology", "2", "05/01/2026", "08:30:00", 720),
  @(74, "Year 4", "INT-B2-B", "dermatology", "3", "06/01/2026", "08:30:00", 720),
  @(75, "Year 4", "INT-B2-B", "dermatology", "4", "07/01/2026", "08:30:00", 720),
  @(76, "Year 4", "INT-B2-B", "dermatology", "5", "08/01/2026", "08:30:00", 720),
  @(77, "Year 4", "INT-B2-B", "immunology/haematology", "1", "23/11/2025", "08:30:00", 720),
  @(78, "Year 4", "INT-B2-B", "immunology/haematology", "2", "24/11/2025", "08:30:00", 720),
  @(79, "Year 4", "INT-B2-B", "immunology/haematology", "3", "25/11/2025", "08:30:00", 720),
  @(80, "Year 4", "INT-B2-B", "immunology/haematology", "4", "26/11/2025", "08:30:00", 720),
  @(81, "Year 4", "INT-B2-B", "immunology/haematology", "5", "27/11/2025", "08:30:00", 720),
  @(82, "Year 4", "INT-B2-B", "immunology/haematology", "6", "30/11/2025", "08:30:00", 720),
  @(83, "Year 4", "INT-B2-B", "immunology/haematology", "7", "01/12/2025", "08:30:00", 720),
  @(84, "Year 4", "INT-B2-B", "immunology/haematology", "8", "02/12/2025", "08:30:00", 720),
  @(85, "Year 4", "INT-B2-B", "immunology/haematology", "9", "03/12/2025", "08:30:00", 720),
  @(86, "Year 4", "INT-B2-B", "immunology/haematology", "10", "04/12/2025", "08:30:00", 720),
  @(87, "Year 4", "INT-B2-B", "tropical", "1", "11/01/2026", "08:30:00", 720),
  @(88, "Year 4", "INT-B2-B", "tropical", "2", "12/01/2026", "08:30:00", 720),
  @(89, "Year 4", "INT-B2-B", "tropical", "3", "13/01/2026", "08:30:00", 720),
  @(90, "Year 4", "INT-B2-B", "tropical", "4", "14/01/2026", "08:30:00", 720),
  @(91, "Year 4", "INT-B2-B", "tropical", "5", "15/01/2026", "08:30:00", 720)
)

$centerH = -4108   # xlCenter
$centerV = -4108   # xlCenter
$bandFillColor = 15790320   # RGB(240,240,240) light-gray row banding, matches existing odd data rows

foreach ($r in $newRows) {
    $rowNum      = $r[0]
    $yearVal     = $r[1]
    $groupVal    = $r[2]
    $subjectVal  = $r[3]
    $sessionVal  = $r[4]
    $dateVal     = $r[5]
    $startTime   = $r[6]
    $duration    = $r[7]

    $isBanded = (($rowNum % 2) -eq 0)

    # Columns A-D: plain text fields (Year, Group, Subject, Session)
    foreach ($pair in @(@("A", $yearVal), @("B", $groupVal), @("C", $subjectVal), @("D", $sessionVal))) {
        $col = $pair[0]
        $val = $pair[1]
        $cell = $ws.Range("$col$rowNum")
        $cell.NumberFormat = "@"
        $cell.HorizontalAlignment = $centerH
        $cell.VerticalAlignment = $centerV
        if ($isBanded) {
            $cell.Interior.Color = $bandFillColor
        }
        $cell.Value = $val
    }

    # Column E: Date, stored as literal text (dd/mm/yyyy display format, text content)
    $eCell = $ws.Range("E$rowNum")
    $eCell.NumberFormat = "@"
    $eCell.HorizontalAlignment = $centerH
    $eCell.VerticalAlignment = $centerV
    if ($isBanded) {
        $eCell.Interior.Color = $bandFillColor
    }
    $eCell.Value = $dateVal

    # Column F: Start time, stored as literal text (hh:mm:ss display format, text content)
    $fCell = $ws.Range("F$rowNum")
    $fCell.NumberFormat = "@"
    $fCell.HorizontalAlignment = $centerH
    $fCell.VerticalAlignment = $centerV
    if ($isBanded) {
        $fCell.Interior.Color = $bandFillColor
    }
    $fCell.Value = $startTime

    # Column G: Duration, numeric (minutes)
    $gCell = $ws.Range("G$rowNum")
    $gCell.NumberFormat = "0"
    $gCell.HorizontalAlignment = $centerH
    $gCell.VerticalAlignment = $centerV
    if ($isBanded) {
        $gCell.Interior.Color = $bandFillColor
    }
    $gCell.Value = $duration
}

Write-Host "Added $($newRows.Count) rows (52-91) for group INT-B2-B."
